$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "-"
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = "[Emerson-Robótica, Euclides-Tecnologia da soldagem, Rogério-Processos de Usinagem 2, Victor Lima-CAM]"
$ws.Range("F2").Value = "-"

# Row 3
$ws.Range("B3").Value = "-"
$ws.Range("C3").Value = "[Pedro Bispo-Lab. Máquinas Elétricas, João Paulo-Lab. de eletroeletrônica, -]"
$ws.Range("D3").Value = "André Guimarães-Máquinas Térmicas e de Fl"
$ws.Range("E3").Value = "[Emerson-Robótica, Euclides-Tecnologia da soldagem, Rogério-Processos de Usinagem 2, Victor Lima-CAM]"
$ws.Range("F3").Value = "[Eudes-Microcontroladores, -, Pedro Bispo-Automação Industrial, Leonardo-Manut. Mecânica]"

# Row 4
$ws.Range("B4").Value = "-"
$ws.Range("C4").Value = "[Pedro Bispo-Lab. Máquinas Elétricas, João Paulo-Lab. de eletroeletrônica, -]"
$ws.Range("D4").Value = "André Guimarães-Máquinas Térmicas e de Fl"
$ws.Range("E4").Value = "[Carlos Eduardo-Processos de Usinagem 1, Nilton Maia-Elementos de máquinas, Humberto-Eletropneumática, Ludoff-Eletrohidráulica]"
$ws.Range("F4").Value = "[Eudes-Microcontroladores, -, Pedro Bispo-Automação Industrial, Leonardo-Manut. Mecânica]"

# Row 6
$ws.Range("B6").Value = "[Victor Lima-CAM, Rogério-Processos de Usinagem 2, Euclides-Tecnologia da soldagem, Emerson-Robótica]"
$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = "Cleidson-Máquinas Elétri"
$ws.Range("E6").Value = "[Carlos Eduardo-Processos de Usinagem 1, Nilton Maia-Elementos de máquinas, Humberto-Eletropneumática, Ludoff-Eletrohidráulica]"
$ws.Range("F6").Value = "[Eudes-Microcontroladores, -, Pedro Bispo-Automação Industrial, Leonardo-Manut. Mecânica]"

# Row 7
$ws.Range("B7").Value = "[Victor Lima-CAM, Rogério-Processos de Usinagem 2, Euclides-Tecnologia da soldagem, Emerson-Robótica]"
$ws.Range("C7").Value = "-"
$ws.Range("D7").Value = "Cleidson-Máquinas Elétri"
$ws.Range("E7").Value = "[Carlos Eduardo-Processos de Usinagem 1, Nilton Maia-Elementos de máquinas, Humberto-Eletropneumática, Ludoff-Eletrohidráulica]"
$ws.Range("F7").Value = "[Eudes-Microcontroladores, -, Pedro Bispo-Automação Industrial, Leonardo-Manut. Mecânica]"

# Row 8
$ws.Range("B8").Value = "-"
$ws.Range("D8").Value = "-"
$ws.Range("E8").Value = "[Carlos Eduardo-Processos de Usinagem 1, Nilton Maia-Elementos de máquinas, Humberto-Eletropneumática, Ludoff-Eletrohidráulica]"
$ws.Range("F8").Value = "-"
